$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.617.45"
$ws.Range("E2").Value = "  +0.90%  "

$ws.Range("D3").Value = "3.155.02"
$ws.Range("E3").Value = "  +0.46%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'609.30"
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("D6").Value = "'144.51"
$ws.Range("E6").Value = "  -1.34%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "3.151.13"
$ws.Range("E8").Value = "  +0.45%  "

$ws.Range("D9").Value = "'0.526"
$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("E10").Value = "  +0.82%  "

$ws.Range("D11").Value = "'5.44"
$ws.Range("E11").Value = "  +2.21%  "

$ws.Range("D12").Value = "'0.470"
$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").Value = "'0.0000259"
$ws.Range("E13").Value = "  +3.79%  "

$ws.Range("D14").Value = "'35.55"
$ws.Range("E14").Value = "  +0.65%  "

$ws.Range("D15").Value = "3.668.48"
$ws.Range("E15").Value = "  +0.35%  "

$ws.Range("E16").Value = "  +2.92%  "

$ws.Range("D17").Value = "64.510.91"
$ws.Range("E17").Value = "  +0.68%  "

$ws.Range("D18").Value = "3.145.90"
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").Value = "'6.89"
$ws.Range("E19").Value = "  +0.39%  "

$ws.Range("D20").Value = "'482.28"
$ws.Range("E20").Value = "  +1.24%  "

$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("D22").Value = "'0.718"
$ws.Range("E22").Value = "  +2.19%  "

$ws.Range("D23").Value = "'7.73"
$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("D24").Value = "'85.45"
$ws.Range("E24").Value = "  +2.41%  "

$ws.Range("D25").Value = "'13.48"
$ws.Range("E25").Value = "  -0.74%  "

$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("E27").Value = "  -1.42%  "

$ws.Range("D28").Value = "'8.48"
$ws.Range("E28").Value = "  +0.97%  "

$ws.Range("D29").Value = "'7.22"
$ws.Range("E29").Value = "  +7.47%  "

$ws.Range("E30").Value = "  +1.93%  "

$ws.Range("D31").Value = "'2.06"
$ws.Range("E31").Value = "  -5.35%  "

$ws.Range("D32").Value = "'27.14"
$ws.Range("E32").Value = "  +3.76%  "

$ws.Range("E33").Value = "  -0.12%  "

$ws.Range("D34").Value = "'2.67"
$ws.Range("E34").Value = "  -2.15%  "

$ws.Range("E35").Value = "  -1.24%  "

$ws.Range("D36").Value = "'6.01"
$ws.Range("E36").Value = "  +0.96%  "

$ws.Range("D37").Value = "0.0₃0772"
$ws.Range("E37").Value = "  +5.57%  "

$ws.Range("D38").Value = "'52.51"
$ws.Range("E38").Value = "  -1.74%  "

$ws.Range("D39").Value = "'3.05"
$ws.Range("E39").Value = "  +5.30%  "

$ws.Range("D40").Value = "'447.58"
$ws.Range("E40").Value = "  -2.90%  "

$ws.Range("D41").Value = "'0.0396"
$ws.Range("E41").Value = "  +0.71%  "

$ws.Range("E42").Value = "  +1.14%  "

$ws.Range("D43").Value = "'8.27"
$ws.Range("E43").Value = "  -1.49%  "

$ws.Range("D44").Value = "2.876.65"
$ws.Range("E44").Value = "  +1.12%  "

$ws.Range("D45").Value = "'0.264"
$ws.Range("E45").Value = "  -0.59%  "

$ws.Range("D46").Value = "'2.25"
$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("D47").Value = "'2.42"
$ws.Range("E47").Value = "  +2.46%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'26.33"
$ws.Range("E48").Value = "  +0.13%  "

$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").Value = "'0.999"
$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("E50").Value = "  -0.28%  "

$ws.Range("D51").Value = "'119.81"
$ws.Range("E51").Value = "  +1.21%  "
